$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Concentrado")
Write-Host ($ws2.Columns.Item(19) | Get-Member -Name "*Width*","*Style*" | Out-String)
